$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(14).Delete()
$ws.Rows.Item(12).Delete()

$ws.Range("E10").ClearContents()
$ws.Range("F10").ClearContents()

$ws.Range("D11").ClearContents()
$ws.Range("F11").ClearContents()
$ws.Range("C11").Value = "…"

$v = $ws.Range("C11").Value()
Write-Host ("C11 now: [{0}]" -f $v)
